# Menus, Updated Icons, Documentation
# In progress changes to menus: simplify the food icons (Beef Burger -> Cheese Burger),
# fill in previously "needed" placeholder rows with real data, and add a bit of styling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 - item renamed from "Beef Burger" to "Cheese Burger"
$ws.Range("A7").Value = "Cheese Burger"

# Row 10 - Hot Dog diet tag updated
$ws.Range("E10").Value = "BC, DF"

# Row 13 - 12 oz Soup: fill in previously blank/"needed" info
$ws.Range("B13").Value = "Item will vary daily - find information in the cafeteria."
$ws.Range("C13").Value = "Item will vary daily - find information in the cafeteria."
$ws.Range("D13").Value = "NA"
$ws.Range("E13").Value = "NA"

# Row 14 - Chef Special: fill in previously blank/"needed" info
$ws.Range("B14").Value = "Item will vary daily - find information in the cafeteria."
$ws.Range("C14").Value = "Item will vary daily - find information in the cafeteria."
$ws.Range("D14").Value = "NA"
$ws.Range("E14").Value = "NA"

# Row 15 - Pepperoni & Sausage Calzone: fill in previously blank/"needed" info
$ws.Range("B15").Value = "Italian Sausage / Pepperoni / Tomato Sauce / Mozzarella / Parmesan"
$ws.Range("C15").Value = "Milk, mustard, wheat. May contain egg, sesame, soy, sulphites."
$ws.Range("D15").Value = "NA"
$ws.Range("E15").Value = "NA"

# Row 16 - 5 Cheese Veggie Calzone: fill in previously blank/"needed" info
$ws.Range("B16").Value = "Spinave / Tomato /Onions / Peppers / Tomato Sauce / Mozzarella / Feta / Parmesan"
$ws.Range("C16").Value = "Milk, mustard, wheat. May contain egg, sesame, soy, sulphites."
$ws.Range("D16").Value = "NA"
$ws.Range("E16").Value = "VEG"

# Row 17 - Pizza Slice: fill in previously blank/"needed" info
$ws.Range("B17").Value = "Item will vary daily - find information in the cafeteria."
$ws.Range("C17").Value = "Soy, wheat, milk, gluten. May contain sulphites."
$ws.Range("D17").Value = "Castle Cheese Mozzarella"
$ws.Range("E17").Value = "BC"

# Update the active selection/cursor position saved with the sheet view
[void]$ws.Range("F32").Select()

# Update workbook window geometry (best effort; some hosts may not persist this)
$win = $excel.ActiveWindow
$win.Left = 0
$win.Top = 0
$win.Width = 960
$win.Height = 1050

$wb.Save()
